$d = $word.ActiveDocument

# map of old -> new text replacements
$replacements = [ordered]@{
    "2023-08-17 Thursday" = "2023-08-18 Friday"
    "95×30=" = "95×81="
    "87×81=" = "37×78="
    "30×29=" = "65×66="
    "89×29=" = "46×40="
    "86×21=" = "15×98="
    "76×55=" = "19×51="
    "72×78=" = "58×53="
    "46×21=" = "64×59="
    "15×69=" = "74×71="
    "76×67=" = "82×45="
    "60×52=" = "33×66="
    "13×98=" = "61×86="
    "77×41=" = "25×13="
    "90×50=" = "80×59="
    "14×17=" = "79×44="
    "93×91=" = "80×86="
    "15×82=" = "56×52="
    "17×15=" = "77×37="
    "34×22=" = "14×92="
    "32×37=" = "87×16="
    "26×74=" = "28×20="
    "23×63=" = "25×55="
    "90×97=" = "19×30="
    "71×44=" = "95×12="
    "49×76=" = "72×85="
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
